# Fruta / hortaliza, semanal
# Add this week's price row for Vega Monumental Concepción - Membrillo,
# pushing the previous week's row down from 18 to 19.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Copy the existing (soon to be superseded) row 18 down to row 19 first,
#    so we don't lose its original values before overwriting row 18.
for ($col = 1; $col -le 20; $col++) {
    $ws.Cells.Item(19, $col).Value2 = $ws.Cells.Item(18, $col).Value2
}
$ws.Range("D19").NumberFormat = $ws.Range("D18").NumberFormat

# 2) Update row 18 with this week's new data (new date, unit-of-sale wording).
$ws.Cells.Item(18, 4).Value2 = 45027
$ws.Cells.Item(18, 17).Value2 = "$/bandeja 18 kilos granel"
$ws.Cells.Item(18, 20).Value2 = 18
